# Updated cryptos list on Tue Nov  7 06:45:49 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns with newly scraped
# figures, and fixes two rows whose rank swapped places (Hedera/WEMIXToken
# at rows 32-33, Gas/MXToken at rows 48-49) by writing the correct
# Coin/Link/Price/Volume for each row.
#
# Price values that look like plain numbers (e.g. "254.04") are written
# with a leading apostrophe so Excel stores them as text (matching the
# original inlineStr cells, e.g. "35.148.75") instead of auto-converting
# them to a numeric type; the cell Style is then reset to "Normal" so the
# quote-prefix doesn't leave a stray number-format behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.143.57'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '1.904.46'
$ws.Range('E4').Value = '  -0.37%  '
$ws.Range('D5').Value = '''254.04'
$ws.Range('E5').Value = '  +3.49%  '
$ws.Range('E6').Value = '  +2.27%  '
$ws.Range('E7').Value = '  -0.41%  '
$ws.Range('D8').Value = '''41.57'
$ws.Range('E8').Value = '  +1.94%  '
$ws.Range('E9').Value = '  +3.06%  '
$ws.Range('D10').Value = '''52.65'
$ws.Range('E10').Value = '  +0.20%  '
$ws.Range('D11').Value = '''0.0760'
$ws.Range('E11').Value = '  +5.87%  '
$ws.Range('D12').Value = '''0.0977'
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('D13').Value = '''13.12'
$ws.Range('E13').Value = '  +3.83%  '
$ws.Range('D14').Value = '2.180.77'
$ws.Range('E14').Value = '  +0.48%  '
$ws.Range('D15').Value = '''0.735'
$ws.Range('E15').Value = '  +4.83%  '
$ws.Range('E16').Value = '  +5.57%  '
$ws.Range('D17').Value = '1.901.63'
$ws.Range('E17').Value = '  +0.37%  '
$ws.Range('D18').Value = '35.136.32'
$ws.Range('D19').Value = '''73.88'
$ws.Range('E19').Value = '  +2.79%  '
$ws.Range('D20').Value = '0.0₃0843'
$ws.Range('E20').Value = '  +3.44%  '
$ws.Range('D21').Value = '''242.88'
$ws.Range('E21').Value = '  +1.15%  '
$ws.Range('D22').Value = '''13.07'
$ws.Range('E22').Value = '  +4.28%  '
$ws.Range('E24').Value = '  -0.43%  '
$ws.Range('D25').Value = '''2.46'
$ws.Range('E25').Value = '  +5.92%  '
$ws.Range('D26').Value = '''2.34'
$ws.Range('E26').Value = '  +0.74%  '
$ws.Range('D27').Value = '''167.75'
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('E28').Value = '  +0.45%  '
$ws.Range('D29').Value = '''18.52'
$ws.Range('E29').Value = '  +1.42%  '
$ws.Range('E30').Value = '  +0.61%  '
$ws.Range('D31').Value = '4.128.20'
$ws.Range('E31').Value = '  -0.35%  '
$ws.Range('B32').Value = 'WEMIXToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D32').Value = '''2.02'
$ws.Range('E32').Value = '  +6.42%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '''0.0600'
$ws.Range('E33').Value = '  +6.11%  '
$ws.Range('E34').Value = '  +4.84%  '
$ws.Range('D35').Value = '''1.61'
$ws.Range('E35').Value = '  +9.34%  '
$ws.Range('E36').Value = '  +4.15%  '
$ws.Range('E37').Value = '  -0.30%  '
$ws.Range('D38').Value = '''0.853'
$ws.Range('E38').Value = '  -7.55%  '
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('D40').Value = '''99.60'
$ws.Range('E40').Value = '  +11.14%  '
$ws.Range('D41').Value = '''17.05'
$ws.Range('E41').Value = '  +4.06%  '
$ws.Range('E42').Value = '  +4.55%  '
$ws.Range('E43').Value = '  +2.01%  '
$ws.Range('D44').Value = '''0.0652'
$ws.Range('E44').Value = '  +2.77%  '
$ws.Range('D45').Value = '''2.42'
$ws.Range('E45').Value = '  +0.77%  '
$ws.Range('D46').Value = '1.310.65'
$ws.Range('E46').Value = '  -2.78%  '
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('B48').Value = 'MXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D48').Value = '''2.75'
$ws.Range('E48').Value = '  -1.18%  '
$ws.Range('B49').Value = 'Gas'
$ws.Range('C49').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D49').Value = '''12.26'
$ws.Range('E49').Value = '  +1.34%  '
$ws.Range('D50').Value = '''6.60'
$ws.Range('E50').Value = '  +2.52%  '
$ws.Range('E51').Value = '  +7.27%  '

# Reset style to default so the forced-text apostrophe prefix does not leave a quotePrefix style
$ws.Range('D5').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
